$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 574; this shifts the existing rows 574-674 down to 575-675
$ws.Rows.Item(574).Insert()

# Populate the newly inserted row 574 with the new record's data
$ws.Cells.Item(574, 1).Value2 = 9
$ws.Cells.Item(574, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(574, 3).Value2 = "Metropolitana"
$ws.Cells.Item(574, 4).Value2 = 44641
$ws.Cells.Item(574, 5).Value2 = 13
$ws.Cells.Item(574, 6).Value2 = "Fruta"
$ws.Cells.Item(574, 7).Value2 = 100102
$ws.Cells.Item(574, 8).Value2 = "Cítricos"
$ws.Cells.Item(574, 9).Value2 = 100102005
$ws.Cells.Item(574, 10).Value2 = "Naranja"
$ws.Cells.Item(574, 11).Value2 = "Valencia"
$ws.Cells.Item(574, 12).Value2 = "Primera"
$ws.Cells.Item(574, 13).Value2 = 470
$ws.Cells.Item(574, 14).Value2 = 10000
$ws.Cells.Item(574, 15).Value2 = 10500
$ws.Cells.Item(574, 16).Value2 = 10234
$ws.Cells.Item(574, 17).Value2 = "`$/caja 15 kilos granel"
$ws.Cells.Item(574, 18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(574, 19).Value2 = 682
$ws.Cells.Item(574, 20).Value2 = 15
